$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append three new result rows (25-27) to the Sonuclar sheet, matching the
# "Dismissal Articles" test case results that were completed.

$ws.Cells.Item(25, 1).Value = "login-functionality;login-with-valid-username-and-password"
$ws.Cells.Item(25, 2).Value = "failed"
$ws.Cells.Item(25, 3).Value = "2021-01-01 17_23_50"
$ws.Cells.Item(25, 4).Value = "chrome"

$ws.Cells.Item(26, 1).Value = "login-functionality;login-with-valid-username-and-password"
$ws.Cells.Item(26, 2).Value = "failed"
$ws.Cells.Item(26, 3).Value = "2021-01-02 19_49_25"
$ws.Cells.Item(26, 4).Value = "chrome"

$ws.Cells.Item(27, 1).Value = "all-steps;step-by-step-from-dismissal-articles"
$ws.Cells.Item(27, 2).Value = "failed"
$ws.Cells.Item(27, 3).Value = "2021-01-02 19_50_39"
$ws.Cells.Item(27, 4).Value = "chrome"
